$wb = $excel.ActiveWorkbook

$wsGlyphs = $wb.Worksheets.Item("Glyphs")
$wsGlyphs.Activate()

$wsGlyphs.Range("A34").Value = "g33"
$wsGlyphs.Range("B34").Value = "Vertical m"

$wsGlyphs.Range("B34").Select()
